$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.236.54"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.241.66"

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  -1.31%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.56"
$ws.Range("E7").Value = "  -3.35%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -3.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.19"
$ws.Range("E10").Value = "  +5.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  -2.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("E13").Value = "  -2.50%  "

$ws.Range("E14").Value = "  -3.36%  "

$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "2.270.96"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").Value = "42.051.44"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").Value = "0.0₃0986"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.02"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("E21").Value = "  +3.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.23"
$ws.Range("E22").Value = "  -1.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.80"
$ws.Range("E23").Value = "  +37.86%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("E26").Value = "  -4.87%  "

$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.00"
$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.61"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0818"
$ws.Range("E31").Value = "  -4.62%  "

$ws.Range("E32").Value = "  -3.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.14"
$ws.Range("E33").Value = "  -3.22%  "

$ws.Range("E34").Value = "  -1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  +11.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.45"
$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0310"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.52"
$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "62.15"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("E42").Value = "  -2.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.56"
$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("E45").Value = "  -2.39%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("E47").Value = "  -2.99%  "

$ws.Range("E48").Value = "  -7.19%  "

$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("E51").Value = "  +0.24%  "
